# Edit: Insert 2 new rows of Cilantro price data at row 327 in the
# "Vega Central Mapocho de Santiago - Cilantro" sheet, shifting the
# existing rows 327:379 down to 329:381.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 327 (existing rows 327-379 shift to 329-381)
$ws.Rows("327:328").Insert()

# --- Fill in the new row 327 ---
$ws.Cells.Item(327, 1).Value = 9
$ws.Cells.Item(327, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(327, 3).Value = "Metropolitana"
$ws.Cells.Item(327, 4).Value = 44505
$ws.Cells.Item(327, 5).Value = 13
$ws.Cells.Item(327, 6).Value = 100112040
$ws.Cells.Item(327, 7).Value = "Cilantro"
$ws.Cells.Item(327, 8).Value = "Sin especificar"
$ws.Cells.Item(327, 9).Value = "Primera"
$ws.Cells.Item(327, 10).Value = 52
$ws.Cells.Item(327, 11).Value = 4000
$ws.Cells.Item(327, 12).Value = 4500
$ws.Cells.Item(327, 13).Value = 4250
$ws.Cells.Item(327, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(327, 15).Value = "Región Metropolitana"
$ws.Cells.Item(327, 16).Value = 118
$ws.Cells.Item(327, 17).Value = 36
$ws.Cells.Item(327, 18).Value = "Hortaliza"

# --- Fill in the new row 328 ---
$ws.Cells.Item(328, 1).Value = 9
$ws.Cells.Item(328, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(328, 3).Value = "Metropolitana"
$ws.Cells.Item(328, 4).Value = 44505
$ws.Cells.Item(328, 5).Value = 13
$ws.Cells.Item(328, 6).Value = 100112040
$ws.Cells.Item(328, 7).Value = "Cilantro"
$ws.Cells.Item(328, 8).Value = "Sin especificar"
$ws.Cells.Item(328, 9).Value = "Primera"
$ws.Cells.Item(328, 10).Value = 160
$ws.Cells.Item(328, 11).Value = 8000
$ws.Cells.Item(328, 12).Value = 10000
$ws.Cells.Item(328, 13).Value = 9000
$ws.Cells.Item(328, 14).Value = "$/docena de atados"
$ws.Cells.Item(328, 15).Value = "Región Metropolitana"
$ws.Cells.Item(328, 16).Value = 3000
$ws.Cells.Item(328, 17).Value = 3
$ws.Cells.Item(328, 18).Value = "Hortaliza"

# Make sure the date column keeps the same date style/format as the rest
# of column D (numFmt for "YYYY-MM-DD HH:MM:SS", same as style used on
# D326 / D329).
$ws.Range("D327:D328").NumberFormat = $ws.Range("D326").NumberFormat
